$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the "RGNR n° <MERGEFIELD $RGNR>" paragraph (not yet known).
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("RGNR n")) {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 2) Remove the "(rif. PM <MERGEFIELD $MAGISTRATO>)" paragraph (not yet
#    known).
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("(rif. PM")) {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 3) Rebalance the "Contravvenzioni" table's 2nd/3rd column widths by
#    1 twip each (4252/4254 -> 4251/4255).
# ---------------------------------------------------------------------
foreach ($tbl in $d.Tables) {
    if ($tbl.Range.Text -like "*Contravvenzioni*") {
        for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
            $tbl.Cell($r, 2).Width = 212.55
            $tbl.Cell($r, 3).Width = 212.75
        }
        break
    }
}

# ---------------------------------------------------------------------
# 4) Give the empty paragraph right after that table an explicit run
#    font (it currently has a bare <w:rPr/>).
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "" -or $p.Range.Text -eq [char]13) {
        $next = $d.Range($p.Range.End, [Math]::Min($p.Range.End + 40, $d.Content.End)).Text
        if ($next -like "Il datore di lavoro*") {
            $r = $p.Range
            $r.Font.NameFarEast = "Noto Serif CJK SC"
            $r.Font.NameBi = "Lohit Devanagari"
            $r.Font.Color = -16777216
            $r.Font.Kerning = 1
            $r.Font.Size = 12
            $r.Font.SizeBi = 12
            $r.LanguageID = "it-IT"
            $r.LanguageIDFarEast = "zh-CN"
            $r.LanguageIDOther = "hi-IN"
            break
        }
    }
}

# ---------------------------------------------------------------------
# 5) Style "Normal": overflowPunct false -> true.
# ---------------------------------------------------------------------
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.HangingPunctuation = 1
